# Applies the "Diverse ändringar av felaktigheter" edits to sjukpenningtal.xlsx
# - Month changes from oktober (10) to november (11) for every data row (B and J columns)
# - A number of sjukpenningtal values in column G are corrected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2..85 (row 1 is the header).
$firstRow = 2
$lastRow = 85

# Update month number (column B) and month name (column J) for every data row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "11"
    $ws.Cells.Item($r, 10).Value = "november"
}

# Updated sjukpenningtal (column G) values, keyed by row number.
$gUpdates = @{
    2  = 10.5
    3  = 13.7
    11 = 10.8
    12 = 14.5
    14 = 11.4
    15 = 15.7
    17 = 14.7
    18 = 19.1
    20 = 18
    21 = 21.1
    22 = 15
    23 = 13.2
    24 = 18
    28 = 0.2
    29 = 6.1
    30 = 8.8
    31 = 3.6
    32 = 14.8
    33 = 21.1
    34 = 8.8
    35 = 15.4
    36 = 21.5
    37 = 9.4
    38 = 17.5
    39 = 23.1
    40 = 12.2
    41 = 19
    42 = 23.2
    44 = 12.7
    45 = 17
    52 = 4.9
    54 = 18.8
    56 = 13.7
    57 = 19.2
    59 = 16.1
    60 = 21.4
    62 = 18.6
    65 = 13
    66 = 17.2
    67 = 9.1
    72 = 7.4
    73 = 4.1
    74 = 13.8
    75 = 18.9
    77 = 14.4
    78 = 19.6
    79 = 9.2
    80 = 17.8
    81 = 23.5
    82 = 12.5
    83 = 20.5
    84 = 23.9
    85 = 17.3
}

foreach ($row in $gUpdates.Keys) {
    $ws.Cells.Item($row, 7).Value = $gUpdates[$row]
}
